$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time-range labels in column C for the first and third pair of
# rows (rows 2-3 and 6-7). Rows 4-5 keep their existing values.
$ws.Range("C2").Value = "2:55-3:0"
$ws.Range("C3").Value = "3:0-3:5"
$ws.Range("C6").Value = "18:55-19:0"
$ws.Range("C7").Value = "19:0-19:5"

# Move the active selection from C11 to B11.
[void]$ws.Range("B11").Select()
